# LidarTotalBOM.xlsx BOM updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "PCB Itself" row: quantity 13 -> 2, description now quotes 5 boards instead of 10
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "4 Layer PCB, price for 5 boards (<50x50mm)"

# "M12 lens adapter" row: quantity 2 -> 1, description now "Price for 1pcs"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "Price for 1pcs"

# "Scanning Mirror" row: quantity 12 -> 3, description now quotes price (12$ for 4 mirrors)
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = "Customized shape, Alum, front surface mirror (12$ for 4 mirrors)"

# Move the visible window/selection up a bit (previously scrolled to A16 / E35 selected)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E27").Select()
